$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 41-45 (columns B through X) with revised figures
$data = @{
    41 = @(12508297, 10364572, 181389, 818830, 61717, 220863, 285759, 575167, 11450719, 3688326, 1266697, 905806, 3215208, 2332903, 41780, 1057577, 1210725, 8066, 689050, 529741, 12516362, 12669510, -153147)
    42 = @(8508745, 6470294, 169949, 763407, 83437, 342879, 149818, 528961, 11453568, 3702381, 1412799, 104817, 4116335, 2083616, 33622, -2944823, 1582943, 3155, 893640, 692459, 8511900, 13039667, -4527767)
    43 = @(10432650, 8473611, 222643, 725814, 94180, 178734, 191345, 546323, 15110236, 3769259, 1458381, 876828, 6817205, 2149949, 38614, -4677586, 1490485, 12378, 829377, 673486, 10445029, 16613100, -6168071)
    44 = @(12814059, 10388695, 445287, 796815, 118475, 148702, 212048, 704036, 13498809, 3921694, 2041135, 51071, 5093571, 2334075, 57262, -684750, 2741679, 4261, 1629821, 1116119, 12818319, 16244748, -3426429)
    45 = @(13409864, 10883348, 253443, 688144, 44917, 103435, 524415, 912162, 13141425, 4011006, 1373007, 916087, 4548705, 2250102, 42519, 268439, 1363169, 3007, 741214, 624962, 13412871, 14507601, -1094730)
    46 = @(14569771, 11730767, 877998, 684890, 87597, 116442, 273560, 798518, 16932658, 4012826, 1592661, 94113, 8904484, 2281739, 46835, -2362888, 1899413, 2198, 1082587, 819025, 14571969, 18834270, -4262301)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# Add new row 46 label (column A) as a new quarterly date entry.
# A leading apostrophe forces Excel to store the value as literal text
# instead of reinterpreting the dd-mm-yyyy-looking string as a date
# serial number; resetting the style afterwards keeps the cell
# unformatted like its neighbours.
$cellA46 = $ws.Cells.Item(46, 1)
$cellA46.Value = "'01-04-2021"
$cellA46.Style = "Normal"

$wb.Save()
